$wb = $excel.ActiveWorkbook

# Resolve worksheets by name so the script is independent of tab ordering.
$ws1 = $wb.Worksheets.Item("展览")      # "Exhibitions" - needs row content updates + trailing row deletion
$ws2 = $wb.Worksheets.Item("演出")      # "Performances" - needs a single numeric update
$ws4 = $wb.Worksheets.Item("全部类型")  # "All types" - needs row content updates + trailing row deletion

# ---------------------------------------------------------------------------
# Sheet "展览": overwrite B2:I20 with the refreshed listing, then drop the
# now-unused trailing rows 21-23 (this also shrinks dimension to A1:I20).
# ---------------------------------------------------------------------------
$c = $ws1.Cells.Item(2,2); $c.NumberFormat = "@"; $c.Value = "2024-05-04"; $c.Style = "Normal"
$ws1.Cells.Item(2,3).Value = "苏州·OCG国潮动漫游戏嘉年华"
$ws1.Cells.Item(2,4).Value = "苏州大道东688号 苏州国际博览中心"
$ws1.Cells.Item(2,5).Value = "2024.05.04 09:00-05.05 17:00"
$ws1.Cells.Item(2,6).Value = 8585
$ws1.Cells.Item(2,7).Value = 75
$ws1.Cells.Item(2,8).Value = "https://show.bilibili.com/platform/detail.html?id=82779"
$ws1.Cells.Item(2,9).Value = "//i1.hdslb.com/bfs/openplatform/202403/hcgdIzw61710298907237.jpeg"
$c = $ws1.Cells.Item(3,2); $c.NumberFormat = "@"; $c.Value = "2024-05-05"; $c.Style = "Normal"
$ws1.Cells.Item(3,3).Value = "张家港·山河月国风动漫节"
$ws1.Cells.Item(3,4).Value = "杨舍镇暨阳中路220号 沙洲宾馆"
$ws1.Cells.Item(3,5).Value = "2024.05.05 11:00-05.05 17:00"
$ws1.Cells.Item(3,6).Value = 72
$ws1.Cells.Item(3,7).Value = 30
$ws1.Cells.Item(3,8).Value = "https://show.bilibili.com/platform/detail.html?id=84432"
$ws1.Cells.Item(3,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/hxPfxsBV1713176346674.jpeg"
$c = $ws1.Cells.Item(4,2); $c.NumberFormat = "@"; $c.Value = "2024-05-18"; $c.Style = "Normal"
$ws1.Cells.Item(4,3).Value = "太仓·原x崩铁ONLY"
$ws1.Cells.Item(4,4).Value = "滨河路128号 凯景世纪大酒店(太仓滨河路店)"
$ws1.Cells.Item(4,5).Value = "2024.05.18 10:00-05.18 17:00"
$ws1.Cells.Item(4,6).Value = 28
$ws1.Cells.Item(4,7).Value = 55
$ws1.Cells.Item(4,8).Value = "https://show.bilibili.com/platform/detail.html?id=84904"
$ws1.Cells.Item(4,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/qxpmdgLs1714118849736.jpeg"
$c = $ws1.Cells.Item(5,2); $c.NumberFormat = "@"; $c.Value = "2024-05-18"; $c.Style = "Normal"
$ws1.Cells.Item(5,3).Value = "苏州·OrangeOrange国潮&随机宅舞派对【免费活动】"
$ws1.Cells.Item(5,4).Value = "狮山路298号 金鹰国际购物中心(狮山路店)"
$ws1.Cells.Item(5,5).Value = "2024.05.18 13:00-05.18 17:00"
$ws1.Cells.Item(5,6).Value = 85
$ws1.Cells.Item(5,7).Value = 29
$ws1.Cells.Item(5,8).Value = "https://show.bilibili.com/platform/detail.html?id=83949"
$ws1.Cells.Item(5,9).Value = "//i1.hdslb.com/bfs/openplatform/202404/DOH6BK8i1712638105049.png"
$c = $ws1.Cells.Item(6,2); $c.NumberFormat = "@"; $c.Value = "2024-05-18"; $c.Style = "Normal"
$ws1.Cells.Item(6,3).Value = "苏州·YoungComic动漫嘉年华"
$ws1.Cells.Item(6,4).Value = "清禾路886号 尹山湖大剧院"
$ws1.Cells.Item(6,5).Value = "2024.05.18 10:00-05.18 17:00"
$ws1.Cells.Item(6,6).Value = 1352
$ws1.Cells.Item(6,7).Value = 60
$ws1.Cells.Item(6,8).Value = "https://show.bilibili.com/platform/detail.html?id=83142"
$ws1.Cells.Item(6,9).Value = "//i2.hdslb.com/bfs/openplatform/202403/4wWLK6Jg1710840463319.jpeg"
$c = $ws1.Cells.Item(7,2); $c.NumberFormat = "@"; $c.Value = "2024-05-19"; $c.Style = "Normal"
$ws1.Cells.Item(7,3).Value = "苏州·国乙only（免费展）"
$ws1.Cells.Item(7,4).Value = "吴中万达广场 吴中万达广场"
$ws1.Cells.Item(7,5).Value = "2024.05.19 14:00-05.19 17:00"
$ws1.Cells.Item(7,6).Value = 118
$ws1.Cells.Item(7,7).Value = 20
$ws1.Cells.Item(7,8).Value = "https://show.bilibili.com/platform/detail.html?id=85022"
$ws1.Cells.Item(7,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/4kJ5GWEo1714137407259.jpeg"
$c = $ws1.Cells.Item(8,2); $c.NumberFormat = "@"; $c.Value = "2024-05-25"; $c.Style = "Normal"
$ws1.Cells.Item(8,3).Value = "苏州·姑苏梦行高校联展"
$ws1.Cells.Item(8,4).Value = "尹山湖商业水街2号楼3层 格莱美婚礼宴会中心(尹山湖商业水街店)"
$ws1.Cells.Item(8,5).Value = "2024.05.25 10:00-05.25 17:30"
$ws1.Cells.Item(8,6).Value = 26
$ws1.Cells.Item(8,7).Value = 60
$ws1.Cells.Item(8,8).Value = "https://show.bilibili.com/platform/detail.html?id=84991"
$ws1.Cells.Item(8,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/ZSJS1FTx1713888009395.png"
$c = $ws1.Cells.Item(9,2); $c.NumberFormat = "@"; $c.Value = "2024-05-25"; $c.Style = "Normal"
$ws1.Cells.Item(9,3).Value = "苏州·燃梦Project"
$ws1.Cells.Item(9,4).Value = "清禾路886号 尹山湖大剧院"
$ws1.Cells.Item(9,5).Value = "2024.05.25 10:30-05.25 16:30"
$ws1.Cells.Item(9,6).Value = 29
$ws1.Cells.Item(9,7).Value = 60
$ws1.Cells.Item(9,8).Value = "https://show.bilibili.com/platform/detail.html?id=83926"
$ws1.Cells.Item(9,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/T3neM9fF1714119047940.jpeg"
$c = $ws1.Cells.Item(10,2); $c.NumberFormat = "@"; $c.Value = "2024-06-08"; $c.Style = "Normal"
$ws1.Cells.Item(10,3).Value = "【会员购严选】苏州·Come in joy动漫国潮文化节"
$ws1.Cells.Item(10,4).Value = "金山南路288号 广电国际会展中心"
$ws1.Cells.Item(10,5).Value = "2024.06.08 10:00-06.09 17:00"
$ws1.Cells.Item(10,6).Value = 9279
$ws1.Cells.Item(10,7).Value = 60
$ws1.Cells.Item(10,8).Value = "https://show.bilibili.com/platform/detail.html?id=82233"
$ws1.Cells.Item(10,9).Value = "//i0.hdslb.com/bfs/openplatform/202403/F86lgbSt1709278264141.jpeg"
$c = $ws1.Cells.Item(11,2); $c.NumberFormat = "@"; $c.Value = "2024-06-15"; $c.Style = "Normal"
$ws1.Cells.Item(11,3).Value = "苏州·蔚蓝档案ONLY#2024~Game Builders Go!!!!"
$ws1.Cells.Item(11,4).Value = "城际路21号 苏州汇融广场假日酒店"
$ws1.Cells.Item(11,5).Value = "2024.06.15 10:00-06.15 17:00"
$ws1.Cells.Item(11,6).Value = 150
$ws1.Cells.Item(11,7).Value = 75
$ws1.Cells.Item(11,8).Value = "https://show.bilibili.com/platform/detail.html?id=84130"
$ws1.Cells.Item(11,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/bpTzFcDq1713253785881.jpeg"
$c = $ws1.Cells.Item(12,2); $c.NumberFormat = "@"; $c.Value = "2024-06-16"; $c.Style = "Normal"
$ws1.Cells.Item(12,3).Value = "苏州·明日方舟ONLY#2024~佑桑柔"
$ws1.Cells.Item(12,4).Value = "城际路21号 苏州汇融广场假日酒店"
$ws1.Cells.Item(12,5).Value = "2024.06.16 10:00-06.16 17:00"
$ws1.Cells.Item(12,6).Value = 94
$ws1.Cells.Item(12,7).Value = 75
$ws1.Cells.Item(12,8).Value = "https://show.bilibili.com/platform/detail.html?id=84046"
$ws1.Cells.Item(12,9).Value = "//i1.hdslb.com/bfs/openplatform/202404/6zZBL5EM1713254151081.jpeg"
$c = $ws1.Cells.Item(13,2); $c.NumberFormat = "@"; $c.Value = "2024-06-29"; $c.Style = "Normal"
$ws1.Cells.Item(13,3).Value = "苏州·归离之缘原神only展"
$ws1.Cells.Item(13,4).Value = "清禾路888号2号楼3楼 格莱美婚礼宴会中心"
$ws1.Cells.Item(13,5).Value = "2024.06.29 09:30-06.29 18:30"
$ws1.Cells.Item(13,6).Value = 217
$ws1.Cells.Item(13,7).Value = 89
$ws1.Cells.Item(13,8).Value = "https://show.bilibili.com/platform/detail.html?id=83271"
$ws1.Cells.Item(13,9).Value = "//i1.hdslb.com/bfs/openplatform/202403/hNkSoRCt1710999968954.png"
$c = $ws1.Cells.Item(14,2); $c.NumberFormat = "@"; $c.Value = "2024-07-06"; $c.Style = "Normal"
$ws1.Cells.Item(14,3).Value = "苏州·第一届寒假动漫展宅舞比赛-CF01"
$ws1.Cells.Item(14,4).Value = "润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店"
$ws1.Cells.Item(14,5).Value = "2024.07.06 10:00-07.06 16:00"
$ws1.Cells.Item(14,6).Value = 173
$ws1.Cells.Item(14,7).Value = 49
$ws1.Cells.Item(14,8).Value = "https://show.bilibili.com/platform/detail.html?id=80528"
$ws1.Cells.Item(14,9).Value = "//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg"
$c = $ws1.Cells.Item(15,2); $c.NumberFormat = "@"; $c.Value = "2024-07-20"; $c.Style = "Normal"
$ws1.Cells.Item(15,3).Value = "苏州·白日梦想7.20全职猎人ONLY展"
$ws1.Cells.Item(15,4).Value = "金芳路与新发路交叉口东南120米 万龙大厦"
$ws1.Cells.Item(15,5).Value = "2024.07.20 09:00-07.20 17:00"
$ws1.Cells.Item(15,6).Value = 349
$ws1.Cells.Item(15,7).Value = 72
$ws1.Cells.Item(15,8).Value = "https://show.bilibili.com/platform/detail.html?id=83508"
$ws1.Cells.Item(15,9).Value = "//i1.hdslb.com/bfs/openplatform/202403/LC3LtiWw1711633827120.jpeg"
$c = $ws1.Cells.Item(16,2); $c.NumberFormat = "@"; $c.Value = "2024-07-20"; $c.Style = "Normal"
$ws1.Cells.Item(16,3).Value = "苏州·萤火国潮文化节动漫品牌博览会"
$ws1.Cells.Item(16,4).Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
$ws1.Cells.Item(16,5).Value = "2024.07.20 10:00-07.21 17:00"
$ws1.Cells.Item(16,6).Value = 6215
$ws1.Cells.Item(16,7).Value = 60
$ws1.Cells.Item(16,8).Value = "https://show.bilibili.com/platform/detail.html?id=83301"
$ws1.Cells.Item(16,9).Value = "//i0.hdslb.com/bfs/openplatform/202403/rV07luU61711274774556.jpeg"
$c = $ws1.Cells.Item(17,2); $c.NumberFormat = "@"; $c.Value = "2024-07-27"; $c.Style = "Normal"
$ws1.Cells.Item(17,3).Value = "苏州·第一届动漫游戏展"
$ws1.Cells.Item(17,4).Value = "清禾路886号 尹山湖大剧院"
$ws1.Cells.Item(17,5).Value = "2024.07.27 10:30-07.27 17:00"
$ws1.Cells.Item(17,6).Value = 1055
$ws1.Cells.Item(17,7).Value = 60
$ws1.Cells.Item(17,8).Value = "https://show.bilibili.com/platform/detail.html?id=84899"
$ws1.Cells.Item(17,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/ARz0BVLv1712661597595.jpeg"
$c = $ws1.Cells.Item(18,2); $c.NumberFormat = "@"; $c.Value = "2024-08-03"; $c.Style = "Normal"
$ws1.Cells.Item(18,3).Value = "苏州·星部落动漫嘉年华"
$ws1.Cells.Item(18,4).Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws1.Cells.Item(18,5).Value = "2024.08.03 09:00-08.04 16:00"
$ws1.Cells.Item(18,6).Value = 72
$ws1.Cells.Item(18,7).Value = 49
$ws1.Cells.Item(18,8).Value = "https://show.bilibili.com/platform/detail.html?id=84858"
$ws1.Cells.Item(18,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/UI5EFZTT1713685680462.jpeg"
$c = $ws1.Cells.Item(19,2); $c.NumberFormat = "@"; $c.Value = "2024-10-01"; $c.Style = "Normal"
$ws1.Cells.Item(19,3).Value = "苏州·第二届Redamancy动漫游戏嘉年华"
$ws1.Cells.Item(19,4).Value = "长江路436号绿宝广场b1 party king运动街区"
$ws1.Cells.Item(19,5).Value = "2024.10.01 10:00-10.05 17:00"
$ws1.Cells.Item(19,6).Value = 40
$ws1.Cells.Item(19,7).Value = 98
$ws1.Cells.Item(19,8).Value = "https://show.bilibili.com/platform/detail.html?id=83576"
$ws1.Cells.Item(19,9).Value = "//i2.hdslb.com/bfs/openplatform/202403/MKyrtd4c1711689984512.jpeg"
$c = $ws1.Cells.Item(20,2); $c.NumberFormat = "@"; $c.Value = "2024-10-01"; $c.Style = "Normal"
$ws1.Cells.Item(20,3).Value = "苏州·第十三届理想乡动漫展-同人创作者大会"
$ws1.Cells.Item(20,4).Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws1.Cells.Item(20,5).Value = "2024.10.01 10:00-10.03 17:00"
$ws1.Cells.Item(20,6).Value = 114
$ws1.Cells.Item(20,7).Value = 39
$ws1.Cells.Item(20,8).Value = "https://show.bilibili.com/platform/detail.html?id=83821"
$ws1.Cells.Item(20,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/OMtuTTFY1711958198579.jpeg"

$ws1.Range("21:23").Delete()

# ---------------------------------------------------------------------------
# Sheet "演出": only the "想去人数" count for row 2 changed.
# ---------------------------------------------------------------------------
$ws2.Cells.Item(2,6).Value = 36

# ---------------------------------------------------------------------------
# Sheet "全部类型": overwrite B2:I22 with the refreshed listing, then drop the
# now-unused trailing rows 23-25 (this also shrinks dimension to A1:I22).
# ---------------------------------------------------------------------------
$c = $ws4.Cells.Item(2,2); $c.NumberFormat = "@"; $c.Value = "2024-05-04"; $c.Style = "Normal"
$ws4.Cells.Item(2,3).Value = "苏州·OCG国潮动漫游戏嘉年华"
$ws4.Cells.Item(2,4).Value = "苏州大道东688号 苏州国际博览中心"
$ws4.Cells.Item(2,5).Value = "2024.05.04 09:00-05.05 17:00"
$ws4.Cells.Item(2,6).Value = 8585
$ws4.Cells.Item(2,7).Value = 75
$ws4.Cells.Item(2,8).Value = "https://show.bilibili.com/platform/detail.html?id=82779"
$ws4.Cells.Item(2,9).Value = "//i1.hdslb.com/bfs/openplatform/202403/hcgdIzw61710298907237.jpeg"
$c = $ws4.Cells.Item(3,2); $c.NumberFormat = "@"; $c.Value = "2024-05-05"; $c.Style = "Normal"
$ws4.Cells.Item(3,3).Value = "张家港·山河月国风动漫节"
$ws4.Cells.Item(3,4).Value = "杨舍镇暨阳中路220号 沙洲宾馆"
$ws4.Cells.Item(3,5).Value = "2024.05.05 11:00-05.05 17:00"
$ws4.Cells.Item(3,6).Value = 72
$ws4.Cells.Item(3,7).Value = 30
$ws4.Cells.Item(3,8).Value = "https://show.bilibili.com/platform/detail.html?id=84432"
$ws4.Cells.Item(3,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/hxPfxsBV1713176346674.jpeg"
$c = $ws4.Cells.Item(4,2); $c.NumberFormat = "@"; $c.Value = "2024-05-18"; $c.Style = "Normal"
$ws4.Cells.Item(4,3).Value = "太仓·原x崩铁ONLY"
$ws4.Cells.Item(4,4).Value = "滨河路128号 凯景世纪大酒店(太仓滨河路店)"
$ws4.Cells.Item(4,5).Value = "2024.05.18 10:00-05.18 17:00"
$ws4.Cells.Item(4,6).Value = 28
$ws4.Cells.Item(4,7).Value = 55
$ws4.Cells.Item(4,8).Value = "https://show.bilibili.com/platform/detail.html?id=84904"
$ws4.Cells.Item(4,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/qxpmdgLs1714118849736.jpeg"
$c = $ws4.Cells.Item(5,2); $c.NumberFormat = "@"; $c.Value = "2024-05-18"; $c.Style = "Normal"
$ws4.Cells.Item(5,3).Value = "苏州·OrangeOrange国潮&随机宅舞派对【免费活动】"
$ws4.Cells.Item(5,4).Value = "狮山路298号 金鹰国际购物中心(狮山路店)"
$ws4.Cells.Item(5,5).Value = "2024.05.18 13:00-05.18 17:00"
$ws4.Cells.Item(5,6).Value = 85
$ws4.Cells.Item(5,7).Value = 29
$ws4.Cells.Item(5,8).Value = "https://show.bilibili.com/platform/detail.html?id=83949"
$ws4.Cells.Item(5,9).Value = "//i1.hdslb.com/bfs/openplatform/202404/DOH6BK8i1712638105049.png"
$c = $ws4.Cells.Item(6,2); $c.NumberFormat = "@"; $c.Value = "2024-05-18"; $c.Style = "Normal"
$ws4.Cells.Item(6,3).Value = "苏州·YoungComic动漫嘉年华"
$ws4.Cells.Item(6,4).Value = "清禾路886号 尹山湖大剧院"
$ws4.Cells.Item(6,5).Value = "2024.05.18 10:00-05.18 17:00"
$ws4.Cells.Item(6,6).Value = 1352
$ws4.Cells.Item(6,7).Value = 60
$ws4.Cells.Item(6,8).Value = "https://show.bilibili.com/platform/detail.html?id=83142"
$ws4.Cells.Item(6,9).Value = "//i2.hdslb.com/bfs/openplatform/202403/4wWLK6Jg1710840463319.jpeg"
$c = $ws4.Cells.Item(7,2); $c.NumberFormat = "@"; $c.Value = "2024-05-19"; $c.Style = "Normal"
$ws4.Cells.Item(7,3).Value = "苏州·国乙only（免费展）"
$ws4.Cells.Item(7,4).Value = "吴中万达广场 吴中万达广场"
$ws4.Cells.Item(7,5).Value = "2024.05.19 14:00-05.19 17:00"
$ws4.Cells.Item(7,6).Value = 118
$ws4.Cells.Item(7,7).Value = 20
$ws4.Cells.Item(7,8).Value = "https://show.bilibili.com/platform/detail.html?id=85022"
$ws4.Cells.Item(7,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/4kJ5GWEo1714137407259.jpeg"
$c = $ws4.Cells.Item(8,2); $c.NumberFormat = "@"; $c.Value = "2024-05-25"; $c.Style = "Normal"
$ws4.Cells.Item(8,3).Value = "苏州·姑苏梦行高校联展"
$ws4.Cells.Item(8,4).Value = "尹山湖商业水街2号楼3层 格莱美婚礼宴会中心(尹山湖商业水街店)"
$ws4.Cells.Item(8,5).Value = "2024.05.25 10:00-05.25 17:30"
$ws4.Cells.Item(8,6).Value = 26
$ws4.Cells.Item(8,7).Value = 60
$ws4.Cells.Item(8,8).Value = "https://show.bilibili.com/platform/detail.html?id=84991"
$ws4.Cells.Item(8,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/ZSJS1FTx1713888009395.png"
$c = $ws4.Cells.Item(9,2); $c.NumberFormat = "@"; $c.Value = "2024-05-25"; $c.Style = "Normal"
$ws4.Cells.Item(9,3).Value = "苏州·燃梦Project"
$ws4.Cells.Item(9,4).Value = "清禾路886号 尹山湖大剧院"
$ws4.Cells.Item(9,5).Value = "2024.05.25 10:30-05.25 16:30"
$ws4.Cells.Item(9,6).Value = 29
$ws4.Cells.Item(9,7).Value = 60
$ws4.Cells.Item(9,8).Value = "https://show.bilibili.com/platform/detail.html?id=83926"
$ws4.Cells.Item(9,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/T3neM9fF1714119047940.jpeg"
$c = $ws4.Cells.Item(10,2); $c.NumberFormat = "@"; $c.Value = "2024-06-01"; $c.Style = "Normal"
$ws4.Cells.Item(10,3).Value = "苏州·春日计划2024——特别二次元不插电音乐会"
$ws4.Cells.Item(10,4).Value = "星湖街555号高教区(体育馆南侧) 苏州独墅湖影剧院"
$ws4.Cells.Item(10,5).Value = "2024.06.01 19:30-06.01 21:00"
$ws4.Cells.Item(10,6).Value = 36
$ws4.Cells.Item(10,7).Value = 88
$ws4.Cells.Item(10,8).Value = "https://show.bilibili.com/platform/detail.html?id=84720"
$ws4.Cells.Item(10,9).Value = "//i1.hdslb.com/bfs/openplatform/202404/gwLWvSew1713796405109.png"
$c = $ws4.Cells.Item(11,2); $c.NumberFormat = "@"; $c.Value = "2024-06-02"; $c.Style = "Normal"
$ws4.Cells.Item(11,3).Value = "苏州·英雄时代2024哈瓦西钢琴演奏会"
$ws4.Cells.Item(11,4).Value = "东太湖大道12000号 苏州湾大剧院"
$ws4.Cells.Item(11,5).Value = "2024.06.02 19:30-06.02 21:00"
$ws4.Cells.Item(11,6).Value = 0
$ws4.Cells.Item(11,7).Value = 499
$ws4.Cells.Item(11,8).Value = "https://show.bilibili.com/platform/detail.html?id=83901"
$ws4.Cells.Item(11,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/LbCirky11712569675168.png"
$c = $ws4.Cells.Item(12,2); $c.NumberFormat = "@"; $c.Value = "2024-06-08"; $c.Style = "Normal"
$ws4.Cells.Item(12,3).Value = "【会员购严选】苏州·Come in joy动漫国潮文化节"
$ws4.Cells.Item(12,4).Value = "金山南路288号 广电国际会展中心"
$ws4.Cells.Item(12,5).Value = "2024.06.08 10:00-06.09 17:00"
$ws4.Cells.Item(12,6).Value = 9279
$ws4.Cells.Item(12,7).Value = 60
$ws4.Cells.Item(12,8).Value = "https://show.bilibili.com/platform/detail.html?id=82233"
$ws4.Cells.Item(12,9).Value = "//i0.hdslb.com/bfs/openplatform/202403/F86lgbSt1709278264141.jpeg"
$c = $ws4.Cells.Item(13,2); $c.NumberFormat = "@"; $c.Value = "2024-06-15"; $c.Style = "Normal"
$ws4.Cells.Item(13,3).Value = "苏州·蔚蓝档案ONLY#2024~Game Builders Go!!!!"
$ws4.Cells.Item(13,4).Value = "城际路21号 苏州汇融广场假日酒店"
$ws4.Cells.Item(13,5).Value = "2024.06.15 10:00-06.15 17:00"
$ws4.Cells.Item(13,6).Value = 150
$ws4.Cells.Item(13,7).Value = 75
$ws4.Cells.Item(13,8).Value = "https://show.bilibili.com/platform/detail.html?id=84130"
$ws4.Cells.Item(13,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/bpTzFcDq1713253785881.jpeg"
$c = $ws4.Cells.Item(14,2); $c.NumberFormat = "@"; $c.Value = "2024-06-16"; $c.Style = "Normal"
$ws4.Cells.Item(14,3).Value = "苏州·明日方舟ONLY#2024~佑桑柔"
$ws4.Cells.Item(14,4).Value = "城际路21号 苏州汇融广场假日酒店"
$ws4.Cells.Item(14,5).Value = "2024.06.16 10:00-06.16 17:00"
$ws4.Cells.Item(14,6).Value = 94
$ws4.Cells.Item(14,7).Value = 75
$ws4.Cells.Item(14,8).Value = "https://show.bilibili.com/platform/detail.html?id=84046"
$ws4.Cells.Item(14,9).Value = "//i1.hdslb.com/bfs/openplatform/202404/6zZBL5EM1713254151081.jpeg"
$c = $ws4.Cells.Item(15,2); $c.NumberFormat = "@"; $c.Value = "2024-06-29"; $c.Style = "Normal"
$ws4.Cells.Item(15,3).Value = "苏州·归离之缘原神only展"
$ws4.Cells.Item(15,4).Value = "清禾路888号2号楼3楼 格莱美婚礼宴会中心"
$ws4.Cells.Item(15,5).Value = "2024.06.29 09:30-06.29 18:30"
$ws4.Cells.Item(15,6).Value = 217
$ws4.Cells.Item(15,7).Value = 89
$ws4.Cells.Item(15,8).Value = "https://show.bilibili.com/platform/detail.html?id=83271"
$ws4.Cells.Item(15,9).Value = "//i1.hdslb.com/bfs/openplatform/202403/hNkSoRCt1710999968954.png"
$c = $ws4.Cells.Item(16,2); $c.NumberFormat = "@"; $c.Value = "2024-07-06"; $c.Style = "Normal"
$ws4.Cells.Item(16,3).Value = "苏州·第一届寒假动漫展宅舞比赛-CF01"
$ws4.Cells.Item(16,4).Value = "润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店"
$ws4.Cells.Item(16,5).Value = "2024.07.06 10:00-07.06 16:00"
$ws4.Cells.Item(16,6).Value = 173
$ws4.Cells.Item(16,7).Value = 49
$ws4.Cells.Item(16,8).Value = "https://show.bilibili.com/platform/detail.html?id=80528"
$ws4.Cells.Item(16,9).Value = "//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg"
$c = $ws4.Cells.Item(17,2); $c.NumberFormat = "@"; $c.Value = "2024-07-20"; $c.Style = "Normal"
$ws4.Cells.Item(17,3).Value = "苏州·白日梦想7.20全职猎人ONLY展"
$ws4.Cells.Item(17,4).Value = "金芳路与新发路交叉口东南120米 万龙大厦"
$ws4.Cells.Item(17,5).Value = "2024.07.20 09:00-07.20 17:00"
$ws4.Cells.Item(17,6).Value = 349
$ws4.Cells.Item(17,7).Value = 72
$ws4.Cells.Item(17,8).Value = "https://show.bilibili.com/platform/detail.html?id=83508"
$ws4.Cells.Item(17,9).Value = "//i1.hdslb.com/bfs/openplatform/202403/LC3LtiWw1711633827120.jpeg"
$c = $ws4.Cells.Item(18,2); $c.NumberFormat = "@"; $c.Value = "2024-07-20"; $c.Style = "Normal"
$ws4.Cells.Item(18,3).Value = "苏州·萤火国潮文化节动漫品牌博览会"
$ws4.Cells.Item(18,4).Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
$ws4.Cells.Item(18,5).Value = "2024.07.20 10:00-07.21 17:00"
$ws4.Cells.Item(18,6).Value = 6215
$ws4.Cells.Item(18,7).Value = 60
$ws4.Cells.Item(18,8).Value = "https://show.bilibili.com/platform/detail.html?id=83301"
$ws4.Cells.Item(18,9).Value = "//i0.hdslb.com/bfs/openplatform/202403/rV07luU61711274774556.jpeg"
$c = $ws4.Cells.Item(19,2); $c.NumberFormat = "@"; $c.Value = "2024-07-27"; $c.Style = "Normal"
$ws4.Cells.Item(19,3).Value = "苏州·第一届动漫游戏展"
$ws4.Cells.Item(19,4).Value = "清禾路886号 尹山湖大剧院"
$ws4.Cells.Item(19,5).Value = "2024.07.27 10:30-07.27 17:00"
$ws4.Cells.Item(19,6).Value = 1055
$ws4.Cells.Item(19,7).Value = 60
$ws4.Cells.Item(19,8).Value = "https://show.bilibili.com/platform/detail.html?id=84899"
$ws4.Cells.Item(19,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/ARz0BVLv1712661597595.jpeg"
$c = $ws4.Cells.Item(20,2); $c.NumberFormat = "@"; $c.Value = "2024-08-03"; $c.Style = "Normal"
$ws4.Cells.Item(20,3).Value = "苏州·星部落动漫嘉年华"
$ws4.Cells.Item(20,4).Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws4.Cells.Item(20,5).Value = "2024.08.03 09:00-08.04 16:00"
$ws4.Cells.Item(20,6).Value = 72
$ws4.Cells.Item(20,7).Value = 49
$ws4.Cells.Item(20,8).Value = "https://show.bilibili.com/platform/detail.html?id=84858"
$ws4.Cells.Item(20,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/UI5EFZTT1713685680462.jpeg"
$c = $ws4.Cells.Item(21,2); $c.NumberFormat = "@"; $c.Value = "2024-10-01"; $c.Style = "Normal"
$ws4.Cells.Item(21,3).Value = "苏州·第二届Redamancy动漫游戏嘉年华"
$ws4.Cells.Item(21,4).Value = "长江路436号绿宝广场b1 party king运动街区"
$ws4.Cells.Item(21,5).Value = "2024.10.01 10:00-10.05 17:00"
$ws4.Cells.Item(21,6).Value = 40
$ws4.Cells.Item(21,7).Value = 98
$ws4.Cells.Item(21,8).Value = "https://show.bilibili.com/platform/detail.html?id=83576"
$ws4.Cells.Item(21,9).Value = "//i2.hdslb.com/bfs/openplatform/202403/MKyrtd4c1711689984512.jpeg"
$c = $ws4.Cells.Item(22,2); $c.NumberFormat = "@"; $c.Value = "2024-10-01"; $c.Style = "Normal"
$ws4.Cells.Item(22,3).Value = "苏州·第十三届理想乡动漫展-同人创作者大会"
$ws4.Cells.Item(22,4).Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws4.Cells.Item(22,5).Value = "2024.10.01 10:00-10.03 17:00"
$ws4.Cells.Item(22,6).Value = 114
$ws4.Cells.Item(22,7).Value = 39
$ws4.Cells.Item(22,8).Value = "https://show.bilibili.com/platform/detail.html?id=83821"
$ws4.Cells.Item(22,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/OMtuTTFY1711958198579.jpeg"

$ws4.Range("23:25").Delete()
